# Add a "Team" column (AF) to the Marcelo Bielsa coaching-career sheet.
# For every match row, AF holds whichever of HomeTeam/AwayTeam corresponds
# to the club Bielsa was managing at the time (his next job always starts
# a fresh contiguous block of rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AF1 = "Team", matching the bold header style used by the
# other header cells (copy formatting from the neighbouring header AE1).
$ws.Range("AE1").Copy() | Out-Null
$ws.Range("AF1").PasteSpecial(-4122) | Out-Null
$ws.Range("AF1").Value2 = "Team"

function Fill-Team($firstRow, $lastRow, $teamName) {
    $count = $lastRow - $firstRow + 1
    $arr = New-Object 'object[,]' $count, 1
    for ($i = 0; $i -lt $count; $i++) {
        $arr[$i, 0] = $teamName
    }
    $rangeAddr = "AF" + $firstRow + ":AF" + $lastRow
    $ws.Range($rangeAddr).Value2 = $arr
}

# Contiguous blocks of rows for each club Bielsa coached during this spell.
Fill-Team 2   77  "Ath Bilbao"
Fill-Team 78  115 "Marseille"
Fill-Team 116 132 "Lille"
Fill-Team 133 260 "Leeds"

$wb.Application.CutCopyMode = $false
